$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # cases_by_race
$ws2 = $wb.Worksheets.Item(2)   # cases_by_ethnicity

# ---------------------------------------------------------------------------
# Sheet 1 (cases_by_race): renumber the running index in column A for the
# existing 2020-12-12 block (rows 35-43) so it continues the running count
# (33-41) instead of restarting at 0, then append the new 2020-12-13 block
# (rows 44-52) with a fresh 0-8 running index.
# ---------------------------------------------------------------------------

$renumber1 = @(
    @(35, 33),
    @(36, 34),
    @(37, 35),
    @(38, 36),
    @(39, 37),
    @(40, 38),
    @(41, 39),
    @(42, 40),
    @(43, 41)
)
foreach ($pair in $renumber1) {
    $ws1.Cells.Item($pair[0], 1).Value = $pair[1]
}

$newRows1 = @(
    @(44, 0, "", "2020-12-13", "2020-12-12", 1),
    @(45, 1, "American Indian or Alaska Native", "2020-12-13", "2020-12-12", 49),
    @(46, 2, "Asian", "2020-12-13", "2020-12-12", 228),
    @(47, 3, "Black or African American", "2020-12-13", "2020-12-12", 1342),
    @(48, 4, "Native Hawaiian or Other Pacific Islander", "2020-12-13", "2020-12-12", 10),
    @(49, 5, "Not disclosed", "2020-12-13", "2020-12-12", 1466),
    @(50, 6, "Other Race", "2020-12-13", "2020-12-12", 356),
    @(51, 7, "Two or more", "2020-12-13", "2020-12-12", 115),
    @(52, 8, "White", "2020-12-13", "2020-12-12", 12356)
)

# Carry the formatting of the last existing data row (column A has the bold
# centered/bordered style, applied via s="1" in the OOXML) down onto the new
# rows before writing values, so the appended block matches the rest of the
# table visually.
$ws1.Range("A43").Copy()
$ws1.Range("A44:A52").PasteSpecial(-4122)

foreach ($row in $newRows1) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]

    if ($row[2] -eq "") {
        $ws1.Cells.Item($r, 2).Formula = "'"
    } else {
        $ws1.Cells.Item($r, 2).Value = $row[2]
    }
    $ws1.Cells.Item($r, 2).Style = "Normal"

    $ws1.Cells.Item($r, 3).NumberFormat = "@"
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 3).Style = "Normal"

    $ws1.Cells.Item($r, 4).NumberFormat = "@"
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 4).Style = "Normal"

    $ws1.Cells.Item($r, 5).Value = $row[5]
}

# ---------------------------------------------------------------------------
# Sheet 2 (cases_by_ethnicity): same pattern - renumber rows 17-19 to
# continue the running index (15-17), then append the new 2020-12-13 block
# (rows 20-22) with a fresh 0-2 running index.
# ---------------------------------------------------------------------------

$renumber2 = @(
    @(17, 15),
    @(18, 16),
    @(19, 17)
)
foreach ($pair in $renumber2) {
    $ws2.Cells.Item($pair[0], 1).Value = $pair[1]
}

$newRows2 = @(
    @(20, 0, "Hispanic or Latino", "2020-12-13", "2020-12-12", 357),
    @(21, 1, "Not Hispanic or Latino", "2020-12-13", "2020-12-12", 12428),
    @(22, 2, "unknown", "2020-12-13", "2020-12-12", 3138)
)

$ws2.Range("A19").Copy()
$ws2.Range("A20:A22").PasteSpecial(-4122)

foreach ($row in $newRows2) {
    $r = $row[0]
    $ws2.Cells.Item($r, 1).Value = $row[1]

    $ws2.Cells.Item($r, 2).Value = $row[2]

    $ws2.Cells.Item($r, 3).NumberFormat = "@"
    $ws2.Cells.Item($r, 3).Value = $row[3]
    $ws2.Cells.Item($r, 3).Style = "Normal"

    $ws2.Cells.Item($r, 4).NumberFormat = "@"
    $ws2.Cells.Item($r, 4).Value = $row[4]
    $ws2.Cells.Item($r, 4).Style = "Normal"

    $ws2.Cells.Item($r, 5).Value = $row[5]
}
